# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.711.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.02%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.789.05"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.17%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.57%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.60%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.787.31"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.15%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("E9").Value = "  +0.02%  "

$ws.Range("E10").Value = "  +0.12%  "

$ws.Range("E11").Value = "  -1.79%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.449"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.32%  "

$ws.Range("E13").Value = "  -1.98%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.97"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.15%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.425.01"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.17%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.791.67"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.95%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.46"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.11%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.703.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.96%  "

$ws.Range("E19").Value = "  +0.92%  "

$ws.Range("E20").Value = "  -0.19%  "

$ws.Range("E21").Value = "  -5.89%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "460.15"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.70%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.696"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.41%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000153"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.42%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.40"
$ws.Range("D25").Style = "Normal"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.02"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.32%  "

$ws.Range("E27").Value = "  -1.80%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.02%  "

$ws.Range("E29").Value = "  +0.17%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.935.50"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.11%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.77"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.73%  "

$ws.Range("E32").Value = "  +4.22%  "

$ws.Range("E33").Value = "  -0.97%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.65"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.46%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.15%  "

$ws.Range("E36").Value = "  -0.21%  "

$ws.Range("E37").Value = "  +0.03%  "

$ws.Range("E38").Value = "  -1.30%  "

$ws.Range("E39").Value = "  +0.88%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.998"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.26%  "

$ws.Range("E41").Value = "  +0.07%  "

$ws.Range("E42").Value = "  -0.06%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "45.92"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.13%  "

$ws.Range("E44").Value = "  +0.00%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "48.10"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.25%  "

$ws.Range("E46").Value = "  -0.43%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "149.43"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.72%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.32"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.58%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "393.64"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.40%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "26.74"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.18%  "

$ws.Range("E51").Value = "  -4.93%  "
